# Update environmental predictors table for manuscript
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (Chlorophyll-a): spatial averaging becomes "Variable*" ---
$ws.Range("F5").Value = "Variable*"

# --- Row 4 (Sea Surface Temperature variance -> standard deviation) ---
$ws.Range("A4").Value = "SST Standard Deviation"

# --- Insert a new row for the "Sex" predictor before the Station row (row 8) ---
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = "Sex"
$ws.Range("B8").Value = "sex"
$ws.Range("C8").Value = "Krill sex"
$ws.Range("D8").Value = "Sexual dimorphism"

# --- Row 4 Purpose column: drop "front probability" wording ---
$ws.Range("D4").Value = "Food availability"

# --- Row 5 (Chlorophyll-a) Temporal Averaging ---
$ws.Range("E5").Value = "27 days"

# --- Row 7 (Coastal Upwelling Transport Index) Temporal Averaging ---
$ws.Range("E7").Value = "9 days"

# --- Fill remaining "Sex" row cells ---
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("G8").Value = "NA"

# Restore the active cell / selection as recorded in the saved workbook
$ws.Range("E7").Select()
